$wb = $excel.ActiveWorkbook

$aboutWs = $wb.Worksheets.Item("About")
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

$aboutWs.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"
$aboutWs.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Abayskaya Coal Mine, Kazakhstan, M1435, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$newVersionText = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"
$dataWs.Range("S2").Value = $newVersionText
$dataWs.Range("S3").Value = $newVersionText
$dataWs.Range("S4").Value = $newVersionText
$dataWs.Range("S5").Value = $newVersionText
$dataWs.Range("S6").Value = $newVersionText
$dataWs.Range("S7").Value = $newVersionText
$dataWs.Range("S8").Value = $newVersionText
$dataWs.Range("S9").Value = $newVersionText
